$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.661.79'
$ws.Range("E2").Value = '  -4.75%  '
$ws.Range("D3").Value = '3.463.91'
$ws.Range("E3").Value = '  -6.27%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.03'
$ws.Range("E5").Value = '  -7.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.53'
$ws.Range("E6").Value = '  -8.33%  '
$ws.Range("D7").Value = '3.466.25'
$ws.Range("E7").Value = '  -6.15%  '
$ws.Range("E9").Value = '  -4.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.138'
$ws.Range("E10").Value = '  -5.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.90'
$ws.Range("E11").Value = '  -4.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.420'
$ws.Range("E12").Value = '  -5.63%  '
$ws.Range("E13").Value = '  -7.04%  '
$ws.Range("D14").Value = '4.050.58'
$ws.Range("E14").Value = '  -6.19%  '
$ws.Range("E15").Value = '  -4.98%  '
$ws.Range("D16").Value = '3.462.63'
$ws.Range("E16").Value = '  -6.14%  '
$ws.Range("D17").Value = '66.588.01'
$ws.Range("E17").Value = '  -4.75%  '
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("E19").Value = '  -2.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.96'
$ws.Range("E20").Value = '  -7.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '441.04'
$ws.Range("E21").Value = '  -6.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.98'
$ws.Range("E22").Value = '  -14.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.622'
$ws.Range("E23").Value = '  -4.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.04'
$ws.Range("E24").Value = '  -3.97%  '
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("D26").Value = '3.602.79'
$ws.Range("E26").Value = '  -6.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000124'
$ws.Range("E27").Value = '  -4.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.04'
$ws.Range("E28").Value = '  -8.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.22'
$ws.Range("E29").Value = '  -10.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.50'
$ws.Range("E30").Value = '  -5.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.57'
$ws.Range("E31").Value = '  -9.03%  '
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("E33").Value = '  -4.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.45'
$ws.Range("E34").Value = '  -5.29%  '
$ws.Range("E35").Value = '  -6.86%  '
$ws.Range("E36").Value = '  -8.36%  '
$ws.Range("D37").Value = '3.450.48'
$ws.Range("E37").Value = '  -6.53%  '
$ws.Range("E38").Value = '  -6.76%  '
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '173.42'
$ws.Range("E41").Value = '  -2.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.15'
$ws.Range("E42").Value = '  -4.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0861'
$ws.Range("E43").Value = '  -5.00%  '
$ws.Range("E44").Value = '  -7.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.879'
$ws.Range("E45").Value = '  -5.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.25'
$ws.Range("E46").Value = '  -3.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.23'
$ws.Range("E47").Value = '  -4.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.06'
$ws.Range("E48").Value = '  -11.79%  '
$ws.Range("E49").Value = '  -13.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.53'
$ws.Range("E50").Value = '  -4.38%  '
$ws.Range("E51").Value = '  -5.16%  '
